$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 5 with the missing Status column value
$ws.Range("F5").Value = "PASS"

# Add new row 6 with a new test account entry
$ws.Range("C6").Value = "vmxtbtigbpuvokh@gmail.com"
$ws.Range("D6").Value = "vbgpuSZVRQ5"
$ws.Range("E6").Value = "pass"
